# Append 12 more rows (206-217) of normalized year/X data to Sheet1,
# continuing the existing A/B series (A: sequential index, B: near-zero
# float residuals), matching the formatting already used for column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$aValues = @(204, 205, 206, 207, 208, 209, 210, 211, 212, 213, 214, 215)
$bValues = @(
    "5.551115123125783E-17",
    "0",
    "4.440892098500626E-17",
    "-3.700743415417188E-17",
    "-2.775557561562891E-17",
    "-6.344131569286608E-17",
    "-3.700743415417188E-17",
    "4.440892098500626E-17",
    "0",
    "0",
    "0",
    "0"
)

$startRow = 206
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = [double]$bValues[$i]
}

$endRow = $startRow + $aValues.Length - 1

# Column A on the existing rows carries a bold/bordered/centered style
# (the same one used for every other row's A cell) - copy that formatting
# onto the newly appended A cells so the look stays consistent.
$ws.Range("A205").Copy() | Out-Null
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
